$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.076"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.68%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08038"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.19%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.37%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.777"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.11%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.043"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9261"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1596"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "47.16%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1908"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08959"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.07%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03456"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.95%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09889"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001397"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.72%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005735"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.69%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.535"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.71%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.77%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.43%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.09%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.037"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.12%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.74%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04503"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.97%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001212"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.18%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004780"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.76%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001228"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.74%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003020"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-32.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.45%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01061"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.36%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007292"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.97%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1335"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002107"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009690"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.32%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006227"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.79%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-63.56%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "10.70%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.14%"
